$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '24.119.90'
$ws.Range("E2").Value = '  -1.30%  '

$ws.Range("D3").Value = '1.640.09'
$ws.Range("E3").Value = '  -1.05%  '

$ws.Range("D4").Value = '0.9975'
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").Value = '311.99'
$ws.Range("E5").Value = '  +0.19%  '

$ws.Range("D6").Value = '0.9977'
$ws.Range("E6").Value = '  -0.27%  '

$ws.Range("D7").Value = '0.3942'
$ws.Range("E7").Value = '  +0.59%  '

$ws.Range("D8").Value = '0.3879'
$ws.Range("E8").Value = '  -1.18%  '

$ws.Range("D9").Value = '52.32'
$ws.Range("E9").Value = '  +4.00%  '

$ws.Range("D10").Value = '1.395'
$ws.Range("E10").Value = '  +0.74%  '

$ws.Range("D11").Value = '0.9976'
$ws.Range("E11").Value = '  -0.27%  '

$ws.Range("D12").Value = '0.08528'
$ws.Range("E12").Value = '  -0.23%  '

$ws.Range("D13").Value = '24.14'
$ws.Range("E13").Value = '  -3.33%  '

$ws.Range("D14").Value = '7.165'
$ws.Range("E14").Value = '  -1.17%  '

$ws.Range("D15").Value = '0.00001310'
$ws.Range("E15").Value = '  +0.31%  '

$ws.Range("D16").Value = '7.699'
$ws.Range("E16").Value = '  +0.78%  '

$ws.Range("D17").Value = '1.636.14'
$ws.Range("E17").Value = '  -1.43%  '

$ws.Range("E18").Value = '  +1.30%  '

$ws.Range("D19").Value = '0.06932'
$ws.Range("E19").Value = '  -0.32%  '

$ws.Range("D20").Value = '20.27'
$ws.Range("E20").Value = '  -3.25%  '

$ws.Range("D21").Value = '6.909'
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = '0.9968'
$ws.Range("E22").Value = '  -0.42%  '

$ws.Range("D23").Value = '13.56'

$ws.Range("D24").Value = '24.102.46'
$ws.Range("E24").Value = '  -1.34%  '

$ws.Range("D25").Value = '2.461'
$ws.Range("E25").Value = '  +4.90%  '

$ws.Range("D26").Value = '2.908'
$ws.Range("E26").Value = '  +4.10%  '

$ws.Range("D27").Value = '22.42'
$ws.Range("E27").Value = '  -1.67%  '

$ws.Range("D28").Value = '157.54'
$ws.Range("E28").Value = '  -1.17%  '

$ws.Range("D29").Value = '141.67'
$ws.Range("E29").Value = '  -2.47%  '

$ws.Range("D30").Value = '5.399'
$ws.Range("E30").Value = '  -5.49%  '

$ws.Range("D31").Value = '8.024'
$ws.Range("E31").Value = '  -2.00%  '

$ws.Range("D32").Value = '2.538'
$ws.Range("E32").Value = '  -3.46%  '

$ws.Range("D33").Value = '1.802.11'
$ws.Range("E33").Value = '  -2.05%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").Value = '0.08218'
$ws.Range("E34").Value = '  -0.37%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '1.016'
$ws.Range("E35").Value = '  -0.49%  '

$ws.Range("D36").Value = '0.02944'
$ws.Range("E36").Value = '  -2.44%  '

$ws.Range("D37").Value = '6.711'
$ws.Range("E37").Value = '  -2.57%  '

$ws.Range("D38").Value = '0.2709'
$ws.Range("E38").Value = '  -2.64%  '

$ws.Range("D39").Value = '0.09228'
$ws.Range("E39").Value = '  -2.70%  '

$ws.Range("D40").Value = '10.46'
$ws.Range("E40").Value = '  +1.92%  '

$ws.Range("D41").Value = '13.86'
$ws.Range("E41").Value = '  +2.98%  '

$ws.Range("D42").Value = '0.7640'
$ws.Range("E42").Value = '  -2.31%  '

$ws.Range("D43").Value = '1.432'
$ws.Range("E43").Value = '  -4.09%  '

$ws.Range("D44").Value = '16.26'
$ws.Range("E44").Value = '  -1.58%  '

$ws.Range("D45").Value = '0.7001'
$ws.Range("E45").Value = '  -0.68%  '

$ws.Range("D46").Value = '2.499'
$ws.Range("E46").Value = '  -2.44%  '

$ws.Range("D47").Value = '4.111'
$ws.Range("E47").Value = '  -1.02%  '

$ws.Range("D48").Value = '0.9974'
$ws.Range("E48").Value = '  -0.28%  '

$ws.Range("D49").Value = '0.08363'
$ws.Range("E49").Value = '  -3.14%  '

$ws.Range("D50").Value = '136.46'
$ws.Range("E50").Value = '  -0.27%  '

$ws.Range("D51").Value = '1.440'
$ws.Range("E51").Value = '  +18.04%  '
